# Update the "想去人数" (want-to-go count) values in column F across all
# four sheets of the workbook, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 926
$ws.Range("F9").Value = 2212
$ws.Range("F10").Value = 631
$ws.Range("F13").Value = 1109
$ws.Range("F14").Value = 187
$ws.Range("F15").Value = 2229
$ws.Range("F16").Value = 685
$ws.Range("F17").Value = 13672
$ws.Range("F19").Value = 1290
$ws.Range("F20").Value = 51
$ws.Range("F21").Value = 562
$ws.Range("F23").Value = 33
$ws.Range("F24").Value = 145
$ws.Range("F25").Value = 81
$ws.Range("F27").Value = 277
$ws.Range("F28").Value = 153

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 18
$ws.Range("F11").Value = 85
$ws.Range("F12").Value = 63
$ws.Range("F17").Value = 22

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5723
$ws.Range("F4").Value = 473

# Sheet 4: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 473
$ws.Range("F8").Value = 926
$ws.Range("F12").Value = 2212
$ws.Range("F13").Value = 631
$ws.Range("F18").Value = 1109
$ws.Range("F20").Value = 187
$ws.Range("F22").Value = 18
$ws.Range("F23").Value = 2229
$ws.Range("F24").Value = 685
$ws.Range("F25").Value = 85
$ws.Range("F26").Value = 63
$ws.Range("F27").Value = 1290
$ws.Range("F28").Value = 51
$ws.Range("F29").Value = 562
$ws.Range("F31").Value = 33
$ws.Range("F32").Value = 145
$ws.Range("F33").Value = 81
$ws.Range("F38").Value = 277
